$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D2:E51 range to text format so numeric-looking strings
# (e.g. "1.007", "0.07152") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.895.97"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.843.89"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "309.54"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "0.4692"
$ws.Range("E7").Value = "  +3.55%  "
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").Value = "0.07152"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "0.9261"
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("D11").Value = "19.58"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "0.07700"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "1.901.70"
$ws.Range("E13").Value = "  +4.85%  "
$ws.Range("D14").Value = "5.287"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "6.399"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "0.000008649"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "26.919.16"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "14.45"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").Value = "5.017"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "1.932"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "151.96"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "18.24"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").Value = "2.009"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "114.20"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "4.879"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").Value = "0.08834"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "3.216"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").Value = "1.179"
$ws.Range("E32").Value = "  +6.22%  "
$ws.Range("D33").Value = "0.7484"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").Value = "2.787"
$ws.Range("E34").Value = "  +1.73%  "
$ws.Range("D35").Value = "4.475"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "1.084"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "0.01939"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "0.05212"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").Value = "2.955"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("D40").Value = "0.5208"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").Value = "6.996"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").Value = "0.1514"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "8.156"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("E44").Value = "  +5.27%  "
$ws.Range("D45").Value = "0.4699"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "1.006"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "101.05"
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("D48").Value = "1.595"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "65.42"
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("D50").Value = "0.06042"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").Value = "0.8944"
$ws.Range("E51").Value = "  +5.55%  "

# Restore the original (unstyled) appearance of the range now that the
# text values have been written, matching the original workbook styling.
$ws.Range("D2:E51").Style = "Normal"

